# Rename the inline picture shapes (Pearson logo x2 in the footers, BTEC
# logo in the header) so their Name no longer collides: the two Pearson
# logos move from "image1.png" to "image2.png" and the BTEC logo moves
# from "image2.jpg" to "image1.jpg".

$d = $word.ActiveDocument
$sec = $d.Sections.First

# --- Header (first-page header holds the BTEC logo) ---------------------
$hdr = $sec.Headers.Item(2)
if ($hdr.Exists) {
    foreach ($shp in $hdr.Range.InlineShapes) {
        if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
            $shp.Name = "image1.jpg"
        }
    }
}

# --- Footers (both the default and first-page footer carry a Pearson
#     logo picture) -------------------------------------------------------
for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        foreach ($shp in $ftr.Range.InlineShapes) {
            if ($shp.AlternativeText -like "*PearsonLogo.png") {
                $shp.Name = "image2.png"
            }
        }
    }
}
